# Guião.docx edit script
# Applies:
#  1. Moves the "_GoBack" bookmark from before "Disponibilidade..." to a new
#     spot in the "Como podemos garantir que a aplicação móvel..." paragraph,
#     right after newly-added "standards " text.
#  2. Rewrites the tail of that same paragraph with new wording about
#     architecture / good-practice standards, replacing the old ending about
#     "vulnerabilidades nos diversos dispositivos...".
#  3. Rewrites the tail of the "cidade inteligente" paragraph, replacing the
#     old "robusto..." ending with a short pointer back to the previous slide.
#  4. Adds a highlighted leading space run before "7º SLIDE".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the _GoBack bookmark from its old location (right before the
#    "Disponibilidade" paragraph's run). We will re-create it further down
#    once the replacement text exists.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 2. "Como podemos garantir que a aplicação móvel..." paragraph: replace the
#    old closing sentence with the new one.
# ---------------------------------------------------------------------------
$old1 = "Para evitar isto tudo, tem de haver um grande sistema de segurança da informação para não haver vulnerabilidades nos diversos dispositivos que estejam conectados à internet e que a partir dele um atacante possa tirar partido."
$new1 = "Para evitar isto tudo, o desenho da arquitetura do sistema deve seguir protocolos standards e quem implementa estes sistemas deve seguir boas práticas para garantir a segurança do sistema em causa."

$r = $d.Content
$found = $r.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# ---------------------------------------------------------------------------
# 3. "Como podemos garantir que uma cidade inteligente..." paragraph: replace
#    the old closing sentence with the new, shorter pointer sentence.
# ---------------------------------------------------------------------------
$old2 = "Para evitar isto também é necessário haver um sistema de segurança de informação robusto e capaz de não ter vulnerabilidades que possam comprometer a vida dos habitantes dessa cidade."
$new2 = "Para tentar evitar isto, a solução é igual à referida no slide anterior."

$r3 = $d.Content
$r3.Find.Execute($old2, $true, $true, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Insert a highlighted leading space before "7º SLIDE".
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("7º SLIDE") | Out-Null
$r4.Collapse(1)
$r4.InsertBefore(" ")

# ---------------------------------------------------------------------------
# 5. Re-anchor the _GoBack bookmark in its new spot, right after
#    "...protocolos standards ". Added last so none of the InsertAfter /
#    InsertBefore calls above can stretch the bookmark's range.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("protocolos standards ") | Out-Null
$r2.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r2)
